# Actualización automática hashcode
# Updates the hashcode values (column B) for the rows identified below,
# matching the new MD5-like hash values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B100" = "aed0b10bec2e9f418199ec1ba1e7d195"
    "B104" = "d3250a5c1d6eed791df60eeb4e0dcd1e"
    "B113" = "d29b6e376f4ab6820f3ee6102e491c52"
    "B122" = "ee278e6bca7a8de6505b4498ce294b23"
    "B164" = "2c5adae7a14dbf122a3e7e333946951c"
    "B230" = "6cc38a03a89a547d65027e64ed10d11e"
    "B233" = "bf8a8bb894e8abfbc617dca6f524925d"
    "B331" = "e40c86b9f34ec9b35c007636ab13d92c"
    "B342" = "b105fc2bcbf5ba38e657e44c9d60c42d"
    "B343" = "9683867abb9ed08c897898b9ce16f688"
    "B619" = "e3ee95ef384d09352f2806899d18ac19"
    "B623" = "ff9f888e91bca8d85efafc7661513a32"
    "B628" = "846c9647ded4ae397a5a92e7ec1d0301"
    "B757" = "4b38c4c301a7b483c17e2da034918036"
    "B760" = "fd27f7cf29ffa4a7b917c79611b0d651"
    "B763" = "e66db9edeb85723f367334b05a32f91c"
    "B767" = "e2db6fe7b898e2dc1173466e6c23a763"
    "B779" = "d7f4356c35eb2b8b0deaac7d4e0be00c"
    "B818" = "1dcbd17e31672161575ab11d3dad0626"
    "B831" = "572bd04638f6b5cead7fe5e2de230d72"
    "B874" = "c9c849f03081bb7a17b5eba5feebb7ea"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
